# Update recomputed TPM-derived metrics for rows 2-7 (Slitrk3-Ptprs LR pair).
# Ligand-expressing-cell count changed 1 -> 2, which ripples through the
# detection rate / expression / specificity / edge-weight columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 5; Value = 2 }
    @{ Row = 2; Col = 6; Value = 1 }
    @{ Row = 2; Col = 7; Value = 0.4189265 }
    @{ Row = 2; Col = 8; Value = 0.837853 }
    @{ Row = 2; Col = 13; Value = 4.296436999999999 }
    @{ Row = 2; Col = 14; Value = 8.592873999999998 }
    @{ Row = 2; Col = 15; Value = 0.08737129157293111 }
    @{ Row = 2; Col = 16; Value = 0.06876644796033347 }
    @{ Row = 2; Col = 17; Value = 1.7998913148805 }
    @{ Row = 2; Col = 18; Value = 7.199565259521998 }
    @{ Row = 2; Col = 19; Value = 0.08737129157293111 }
    @{ Row = 2; Col = 20; Value = 0.06876644796033347 }
    @{ Row = 3; Col = 5; Value = 2 }
    @{ Row = 3; Col = 6; Value = 1 }
    @{ Row = 3; Col = 7; Value = 0.4189265 }
    @{ Row = 3; Col = 8; Value = 0.837853 }
    @{ Row = 3; Col = 15; Value = 0.2709785829485105 }
    @{ Row = 3; Col = 16; Value = 0.3199146015909443 }
    @{ Row = 3; Col = 17; Value = 5.582291267384166 }
    @{ Row = 3; Col = 18; Value = 33.493747604305 }
    @{ Row = 3; Col = 19; Value = 0.2709785829485105 }
    @{ Row = 3; Col = 20; Value = 0.3199146015909443 }
    @{ Row = 4; Col = 5; Value = 2 }
    @{ Row = 4; Col = 6; Value = 1 }
    @{ Row = 4; Col = 7; Value = 0.4189265 }
    @{ Row = 4; Col = 8; Value = 0.837853 }
    @{ Row = 4; Col = 13; Value = 6.89049 }
    @{ Row = 4; Col = 14; Value = 20.67147 }
    @{ Row = 4; Col = 15; Value = 0.140123318663899 }
    @{ Row = 4; Col = 16; Value = 0.1654281868928364 }
    @{ Row = 4; Col = 17; Value = 2.886608858985 }
    @{ Row = 4; Col = 18; Value = 17.31965315391 }
    @{ Row = 4; Col = 19; Value = 0.140123318663899 }
    @{ Row = 4; Col = 20; Value = 0.1654281868928364 }
    @{ Row = 5; Col = 5; Value = 2 }
    @{ Row = 5; Col = 6; Value = 1 }
    @{ Row = 5; Col = 7; Value = 0.4189265 }
    @{ Row = 5; Col = 8; Value = 0.837853 }
    @{ Row = 5; Col = 13; Value = 18.2696115 }
    @{ Row = 5; Col = 14; Value = 36.539223 }
    @{ Row = 5; Col = 15; Value = 0.3715263492262718 }
    @{ Row = 5; Col = 16; Value = 0.292413525083752 }
    @{ Row = 5; Col = 17; Value = 7.653624402054749 }
    @{ Row = 5; Col = 18; Value = 30.614497608219 }
    @{ Row = 5; Col = 19; Value = 0.3715263492262718 }
    @{ Row = 5; Col = 20; Value = 0.292413525083752 }
    @{ Row = 6; Col = 5; Value = 2 }
    @{ Row = 6; Col = 6; Value = 1 }
    @{ Row = 6; Col = 7; Value = 0.4189265 }
    @{ Row = 6; Col = 8; Value = 0.837853 }
    @{ Row = 6; Col = 13; Value = 1.355562 }
    @{ Row = 6; Col = 14; Value = 4.066686000000001 }
    @{ Row = 6; Col = 15; Value = 0.02756637715092428 }
    @{ Row = 6; Col = 16; Value = 0.03254458882907125 }
    @{ Row = 6; Col = 17; Value = 0.567880844193 }
    @{ Row = 6; Col = 18; Value = 3.407285065158 }
    @{ Row = 6; Col = 19; Value = 0.02756637715092428 }
    @{ Row = 6; Col = 20; Value = 0.03254458882907125 }
    @{ Row = 7; Col = 5; Value = 2 }
    @{ Row = 7; Col = 6; Value = 1 }
    @{ Row = 7; Col = 7; Value = 0.4189265 }
    @{ Row = 7; Col = 8; Value = 0.837853 }
    @{ Row = 7; Col = 13; Value = 5.037141666666667 }
    @{ Row = 7; Col = 14; Value = 15.111425 }
    @{ Row = 7; Col = 15; Value = 0.1024340804374633 }
    @{ Row = 7; Col = 16; Value = 0.1209326496430627 }
    @{ Row = 7; Col = 17; Value = 2.110192128420833 }
    @{ Row = 7; Col = 18; Value = 12.661152770525 }
    @{ Row = 7; Col = 19; Value = 0.1024340804374633 }
    @{ Row = 7; Col = 20; Value = 0.1209326496430627 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}
